# "Fixing sorting in part 2 graphs"
#
# The working-set-size column (C2:C19) is formatted as a number with two
# decimal places so the line chart's category axis sorts/labels the values
# numerically instead of as "General" text-ish values. This also updates the
# chart's category-axis number format to match, and moves the active
# selection/scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply a "0.00" number format to the working-set-size values used as the
# chart's category axis (C2:C19). This creates a new cell style (xfId=1,
# numFmtId=2 "0.00") and stamps it onto every cell in the range.
$ws.Range("C2:C19").NumberFormat = "0.00"

# Match the chart's category axis number format/labels to the new "0.00"
# format on the source data.
$co = $ws.ChartObjects(1)
$chart = $co.Chart
$catAx = $chart.Axes(1)
$catAx.TickLabels.NumberFormat = "0.00"

# Move the selection/view down to F20 (previously F14, scrolled to A3).
$ws.Range("F20").Select() | Out-Null
